$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6089751124382019
$ws.Range("B1").Value = 1.532878279685974
$ws.Range("C1").Value = 5.282925605773926
$ws.Range("D1").Value = 3.046599149703979
$ws.Range("E1").Value = 0.7853374481201172
